$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-07 Tuesday" "2025-01-08 Wednesday"

Replace-Text "67×33=2211" "78×51=3978"
Replace-Text "44×43=1892" "59×11=649"
Replace-Text "21×84=1764" "22×48=1056"
Replace-Text "33×71=2343" "77×62=4774"
Replace-Text "48×94=4512" "78×38=2964"

Replace-Text "52×94=4888" "37×48=1776"
Replace-Text "75×50=3750" "85×50=4250"
Replace-Text "15×64=960" "79×71=5609"
Replace-Text "72×68=4896" "16×86=1376"
Replace-Text "45×54=2430" "81×12=972"

Replace-Text "54×24=1296" "69×73=5037"
Replace-Text "55×41=2255" "96×98=9408"
Replace-Text "38×65=2470" "64×95=6080"
Replace-Text "55×34=1870" "89×44=3916"
Replace-Text "88×88=7744" "36×83=2988"

Replace-Text "52×27=1404" "47×42=1974"
Replace-Text "65×74=4810" "29×36=1044"
Replace-Text "85×67=5695" "21×31=651"
Replace-Text "48×25=1200" "22×65=1430"
Replace-Text "47×81=3807" "22×88=1936"

Replace-Text "35×25=875" "12×87=1044"
Replace-Text "90×81=7290" "31×85=2635"
Replace-Text "30×17=510" "73×56=4088"
Replace-Text "67×38=2546" "46×49=2254"
Replace-Text "86×93=7998" "44×89=3916"
